# Auto-generated Excel COM-interop script to update Typhon_Profits leve-price data
# Applies per-cell numeric updates (and a couple of cell clears) across 8 worksheets
# as captured by the upstream scheduled-runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2374.25  # H62
$ws.Cells.Item(62, 9).Value = 2374.25  # I62
$ws.Cells.Item(62, 11).Value = 2374.25  # K62
$ws.Cells.Item(62, 13).Value = -1750.25  # M62

$ws.Cells.Item(65, 8).Value = 2374.25  # H65
$ws.Cells.Item(65, 9).Value = 2374.25  # I65
$ws.Cells.Item(65, 11).Value = 11871.25  # K65
$ws.Cells.Item(65, 13).Value = -8751.25  # M65

$ws.Cells.Item(86, 8).Value = 9006.429  # H86
$ws.Cells.Item(86, 9).Value = 2220  # I86
$ws.Cells.Item(86, 10).Value = 12776.667  # J86
$ws.Cells.Item(86, 11).Value = 2220  # K86
$ws.Cells.Item(86, 12).Value = 12776.667  # L86
$ws.Cells.Item(86, 13).Value = -1097  # M86
$ws.Cells.Item(86, 14).Value = -15022.667  # N86

$ws.Cells.Item(89, 8).Value = 9006.429  # H89
$ws.Cells.Item(89, 9).Value = 2220  # I89
$ws.Cells.Item(89, 10).Value = 12776.667  # J89
$ws.Cells.Item(89, 11).Value = 11100  # K89
$ws.Cells.Item(89, 12).Value = 63883.335  # L89
$ws.Cells.Item(89, 13).Value = -5484  # M89
$ws.Cells.Item(89, 14).Value = -75115.33499999999  # N89

$ws.Cells.Item(98, 8).Value = 310.44446  # H98
$ws.Cells.Item(98, 9).Value = 310.44446  # I98
$ws.Cells.Item(98, 11).Value = 310.44446  # K98
$ws.Cells.Item(98, 13).Value = 1187.55554  # M98

$ws.Cells.Item(122, 8).Value = 310.44446  # H122
$ws.Cells.Item(122, 9).Value = 310.44446  # I122
$ws.Cells.Item(122, 11).Value = 931.33338  # K122
$ws.Cells.Item(122, 13).Value = 1518.66662  # M122

$ws.Cells.Item(131, 8).Value = 3897.5  # H131
$ws.Cells.Item(131, 9).Value = 0  # I131
$ws.Cells.Item(131, 10).Value = 3897.5  # J131
$ws.Cells.Item(131, 11).Value = 0  # K131
$ws.Cells.Item(131, 12).Value = 11692.5  # L131
$ws.Cells.Item(131, 13).ClearContents()  # M131
$ws.Cells.Item(131, 14).Value = -21772.5  # N131

$ws.Cells.Item(132, 8).Value = 30865.428  # H132
$ws.Cells.Item(132, 9).Value = 33726.25  # I132
$ws.Cells.Item(132, 10).Value = 350  # J132
$ws.Cells.Item(132, 11).Value = 101178.75  # K132
$ws.Cells.Item(132, 12).Value = 1050  # L132
$ws.Cells.Item(132, 13).Value = -98648.75  # M132
$ws.Cells.Item(132, 14).Value = -6110  # N132

$ws.Cells.Item(135, 8).Value = 29416038  # H135
$ws.Cells.Item(135, 9).Value = 1333.091  # I135
$ws.Cells.Item(135, 11).Value = 11997.819  # K135
$ws.Cells.Item(135, 13).Value = -9462.819  # M135

$ws.Cells.Item(139, 8).Value = 52630  # H139
$ws.Cells.Item(139, 10).Value = 52630  # J139
$ws.Cells.Item(139, 12).Value = 52630  # L139
$ws.Cells.Item(139, 14).Value = -62910  # N139

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 350.66666  # H4
$ws.Cells.Item(4, 9).Value = 350  # I4
$ws.Cells.Item(4, 10).Value = 351  # J4
$ws.Cells.Item(4, 11).Value = 350  # K4
$ws.Cells.Item(4, 12).Value = 351  # L4
$ws.Cells.Item(4, 13).Value = -234  # M4
$ws.Cells.Item(4, 14).Value = -583  # N4

$ws.Cells.Item(32, 8).Value = 2883.375  # H32
$ws.Cells.Item(32, 9).Value = 1774.0857  # I32
$ws.Cells.Item(32, 10).Value = 5869.923  # J32
$ws.Cells.Item(32, 11).Value = 1774.0857  # K32
$ws.Cells.Item(32, 12).Value = 5869.923  # L32
$ws.Cells.Item(32, 13).Value = -1487.0857  # M32
$ws.Cells.Item(32, 14).Value = -6443.923  # N32

$ws.Cells.Item(45, 8).Value = 2070.0386  # H45
$ws.Cells.Item(45, 9).Value = 1499.8422  # I45
$ws.Cells.Item(45, 10).Value = 3617.7144  # J45
$ws.Cells.Item(45, 11).Value = 1499.8422  # K45
$ws.Cells.Item(45, 12).Value = 3617.7144  # L45
$ws.Cells.Item(45, 13).Value = -1122.8422  # M45
$ws.Cells.Item(45, 14).Value = -4371.7144  # N45

$ws.Cells.Item(102, 8).Value = 3071.2354  # H102
$ws.Cells.Item(102, 9).Value = 3366.6667  # I102
$ws.Cells.Item(102, 10).Value = 3007.9285  # J102
$ws.Cells.Item(102, 11).Value = 3366.6667  # K102
$ws.Cells.Item(102, 12).Value = 3007.9285  # L102
$ws.Cells.Item(102, 13).Value = -1744.6667  # M102
$ws.Cells.Item(102, 14).Value = -6251.9285  # N102

$ws.Cells.Item(132, 8).Value = 28897.316  # H132
$ws.Cells.Item(132, 9).Value = 1892.3  # I132
$ws.Cells.Item(132, 11).Value = 5676.9  # K132
$ws.Cells.Item(132, 13).Value = -3146.9  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 3226.652  # H94
$ws.Cells.Item(94, 9).Value = 1119.4  # I94
$ws.Cells.Item(94, 11).Value = 1119.4  # K94
$ws.Cells.Item(94, 13).Value = -668.4000000000001  # M94

$ws.Cells.Item(105, 8).Value = 2942973.8  # H105
$ws.Cells.Item(105, 9).Value = 1772.6666  # I105
$ws.Cells.Item(105, 10).Value = 6251825  # J105
$ws.Cells.Item(105, 11).Value = 1772.6666  # K105
$ws.Cells.Item(105, 12).Value = 6251825  # L105
$ws.Cells.Item(105, 13).Value = -25.66660000000002  # M105
$ws.Cells.Item(105, 14).Value = -6255319  # N105

$ws.Cells.Item(107, 8).Value = 1426.375  # H107
$ws.Cells.Item(107, 9).Value = 1282.2  # I107
$ws.Cells.Item(107, 11).Value = 1282.2  # K107
$ws.Cells.Item(107, 13).Value = 637.8  # M107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9489.77  # H31
$ws.Cells.Item(31, 9).Value = 10735.7  # I31
$ws.Cells.Item(31, 10).Value = 5336.6665  # J31
$ws.Cells.Item(31, 11).Value = 10735.7  # K31
$ws.Cells.Item(31, 12).Value = 5336.6665  # L31
$ws.Cells.Item(31, 13).Value = -10440.7  # M31
$ws.Cells.Item(31, 14).Value = -5926.6665  # N31

$ws.Cells.Item(34, 8).Value = 9489.77  # H34
$ws.Cells.Item(34, 9).Value = 10735.7  # I34
$ws.Cells.Item(34, 10).Value = 5336.6665  # J34
$ws.Cells.Item(34, 11).Value = 10735.7  # K34
$ws.Cells.Item(34, 12).Value = 5336.6665  # L34
$ws.Cells.Item(34, 13).Value = -10533.7  # M34
$ws.Cells.Item(34, 14).Value = -5740.6665  # N34

$ws.Cells.Item(107, 8).Value = 689.5  # H107
$ws.Cells.Item(107, 9).Value = 689.5  # I107
$ws.Cells.Item(107, 10).Value = 0  # J107
$ws.Cells.Item(107, 11).Value = 689.5  # K107
$ws.Cells.Item(107, 12).Value = 0  # L107
$ws.Cells.Item(107, 13).Value = 1230.5  # M107
$ws.Cells.Item(107, 14).ClearContents()  # N107

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 10000067  # H4
$ws.Cells.Item(4, 9).Value = 200  # I4
$ws.Cells.Item(4, 11).Value = 600  # K4
$ws.Cells.Item(4, 13).Value = -488  # M4

$ws.Cells.Item(92, 8).Value = 1160  # H92
$ws.Cells.Item(92, 9).Value = 1000  # I92
$ws.Cells.Item(92, 10).Value = 1200  # J92
$ws.Cells.Item(92, 11).Value = 3000  # K92
$ws.Cells.Item(92, 12).Value = 3600  # L92
$ws.Cells.Item(92, 13).Value = -1752  # M92
$ws.Cells.Item(92, 14).Value = -6096  # N92

$ws.Cells.Item(97, 8).Value = 942.3  # H97
$ws.Cells.Item(97, 10).Value = 1179.7142  # J97
$ws.Cells.Item(97, 12).Value = 3539.1426  # L97
$ws.Cells.Item(97, 14).Value = -4531.142599999999  # N97

$ws.Cells.Item(131, 8).Value = 814.1900000000001  # H131
$ws.Cells.Item(131, 10).Value = 826.5876500000001  # J131
$ws.Cells.Item(131, 12).Value = 2479.76295  # L131
$ws.Cells.Item(131, 14).Value = -12559.76295  # N131

$ws.Cells.Item(140, 8).Value = 5194.5806  # H140
$ws.Cells.Item(140, 9).Value = 7296.1875  # I140
$ws.Cells.Item(140, 10).Value = 2952.8667  # J140
$ws.Cells.Item(140, 11).Value = 21888.5625  # K140
$ws.Cells.Item(140, 12).Value = 8858.6001  # L140
$ws.Cells.Item(140, 13).Value = -16708.5625  # M140
$ws.Cells.Item(140, 14).Value = -19218.6001  # N140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3394.7222  # H80
$ws.Cells.Item(80, 10).Value = 4600  # J80
$ws.Cells.Item(80, 12).Value = 4600  # L80
$ws.Cells.Item(80, 14).Value = -6596  # N80

$ws.Cells.Item(83, 8).Value = 3394.7222  # H83
$ws.Cells.Item(83, 10).Value = 4600  # J83
$ws.Cells.Item(83, 12).Value = 23000  # L83
$ws.Cells.Item(83, 14).Value = -32984  # N83

$ws.Cells.Item(97, 8).Value = 1639.5186  # H97
$ws.Cells.Item(97, 9).Value = 887.1579  # I97
$ws.Cells.Item(97, 11).Value = 887.1579  # K97
$ws.Cells.Item(97, 13).Value = -391.1579  # M97

$ws.Cells.Item(102, 8).Value = 2294.4644  # H102
$ws.Cells.Item(102, 9).Value = 2517.9546  # I102
$ws.Cells.Item(102, 11).Value = 2517.9546  # K102
$ws.Cells.Item(102, 13).Value = -895.9546  # M102

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3491.08  # H7
$ws.Cells.Item(7, 9).Value = 3741.6667  # I7
$ws.Cells.Item(7, 10).Value = 2846.7144  # J7
$ws.Cells.Item(7, 11).Value = 3741.6667  # K7
$ws.Cells.Item(7, 12).Value = 2846.7144  # L7
$ws.Cells.Item(7, 13).Value = -3629.6667  # M7
$ws.Cells.Item(7, 14).Value = -3070.7144  # N7

$ws.Cells.Item(93, 8).Value = 2889.95  # H93
$ws.Cells.Item(93, 9).Value = 2708.75  # I93
$ws.Cells.Item(93, 10).Value = 3161.75  # J93
$ws.Cells.Item(93, 11).Value = 2708.75  # K93
$ws.Cells.Item(93, 12).Value = 3161.75  # L93
$ws.Cells.Item(93, 13).Value = -1460.75  # M93
$ws.Cells.Item(93, 14).Value = -5657.75  # N93

$ws.Cells.Item(126, 8).Value = 3491.08  # H126
$ws.Cells.Item(126, 9).Value = 3741.6667  # I126
$ws.Cells.Item(126, 10).Value = 2846.7144  # J126
$ws.Cells.Item(126, 11).Value = 11225.0001  # K126
$ws.Cells.Item(126, 12).Value = 8540.143199999999  # L126
$ws.Cells.Item(126, 13).Value = -8755.000100000001  # M126
$ws.Cells.Item(126, 14).Value = -13480.1432  # N126

$ws.Cells.Item(132, 8).Value = 1583.8928  # H132
$ws.Cells.Item(132, 9).Value = 1049.5264  # I132
$ws.Cells.Item(132, 11).Value = 3148.5792  # K132
$ws.Cells.Item(132, 13).Value = -618.5792000000001  # M132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1601.5946  # H132
$ws.Cells.Item(132, 9).Value = 1486.409  # I132
$ws.Cells.Item(132, 10).Value = 1770.5333  # J132
$ws.Cells.Item(132, 11).Value = 4459.227000000001  # K132
$ws.Cells.Item(132, 12).Value = 5311.5999  # L132
$ws.Cells.Item(132, 13).Value = -1929.227000000001  # M132
$ws.Cells.Item(132, 14).Value = -10371.5999  # N132
